$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2339181286549707
$ws.Range("C2").Value = 0.5058479532163743
$ws.Range("J2").Value = 0.01754385964912281
$ws.Range("P2").Value = 0.1549707602339181
$ws.Range("S2").Value = 0.08771929824561403
$ws.Range("B3").Value = 0.01136363636363636
$ws.Range("C3").Value = 0.005681818181818182
$ws.Range("J3").Value = 0.01704545454545454
$ws.Range("P3").Value = 0.7215909090909091
$ws.Range("S3").Value = 0.2443181818181818
$ws.Range("J4").Value = 0.07407407407407407
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2592592592592592
$ws.Range("B6").Value = 0.06111111111111111
$ws.Range("D6").Value = 0.01666666666666667
$ws.Range("F6").Value = 0.08333333333333333
$ws.Range("J6").Value = 0.2888888888888889
$ws.Range("O6").Value = 0.01111111111111111
$ws.Range("Q6").Value = 0.1944444444444444
$ws.Range("R6").Value = 0.05555555555555555
$ws.Range("S6").Value = 0.2888888888888889
$ws.Range("B7").Value = 0.1007194244604317
$ws.Range("F7").Value = 0.07194244604316546
$ws.Range("J7").Value = 0.1726618705035971
$ws.Range("O7").Value = 0.01438848920863309
$ws.Range("Q7").Value = 0.1223021582733813
$ws.Range("R7").Value = 0.06474820143884892
$ws.Range("S7").Value = 0.4532374100719425
$ws.Range("B8").Value = 0.1355140186915888
$ws.Range("D8").Value = 0.03271028037383177
$ws.Range("E8").Value = 0.004672897196261682
$ws.Range("F8").Value = 0.04205607476635514
$ws.Range("J8").Value = 0.1378504672897196
$ws.Range("O8").Value = 0.02570093457943925
$ws.Range("Q8").Value = 0.161214953271028
$ws.Range("R8").Value = 0.06775700934579439
$ws.Range("S8").Value = 0.3925233644859813
$ws.Range("B9").Value = 0.0684931506849315
$ws.Range("D9").Value = 0.0273972602739726
$ws.Range("F9").Value = 0.02054794520547945
$ws.Range("J9").Value = 0.1712328767123288
$ws.Range("Q9").Value = 0.2465753424657534
$ws.Range("R9").Value = 0.0958904109589041
$ws.Range("S9").Value = 0.3698630136986301
$ws.Range("B10").Value = 0.1315577078288943
$ws.Range("D10").Value = 0.0258272800645682
$ws.Range("E10").Value = 0.002421307506053269
$ws.Range("F10").Value = 0.05811138014527845
$ws.Range("J10").Value = 0.1146085552865214
$ws.Range("O10").Value = 0.01937046004842615
$ws.Range("Q10").Value = 0.2227602905569007
$ws.Range("R10").Value = 0.09362389023405973
$ws.Range("S10").Value = 0.3317191283292978
$ws.Range("G11").Value = 0.1511111111111111
$ws.Range("J11").Value = 0.08888888888888889
$ws.Range("K11").Value = 0.2044444444444445
$ws.Range("L11").Value = 0.5422222222222223
$ws.Range("S11").Value = 0.01333333333333333
$ws.Range("G12").Value = 0.7258064516129032
$ws.Range("J12").Value = 0.2016129032258064
$ws.Range("L12").Value = 0.01612903225806452
$ws.Range("S12").Value = 0.0564516129032258
$ws.Range("G13").Value = 0.6470588235294118
$ws.Range("J13").Value = 0.2941176470588235
$ws.Range("S13").Value = 0.05882352941176471
$ws.Range("F15").Value = 0.02072538860103627
$ws.Range("H15").Value = 0.1398963730569948
$ws.Range("I15").Value = 0.04145077720207254
$ws.Range("J15").Value = 0.4352331606217616
$ws.Range("K15").Value = 0.05699481865284974
$ws.Range("M15").Value = 0.02072538860103627
$ws.Range("O15").Value = 0.02072538860103627
$ws.Range("S15").Value = 0.2642487046632124
$ws.Range("F16").Value = 0.01904761904761905
$ws.Range("H16").Value = 0.1666666666666667
$ws.Range("I16").Value = 0.05714285714285714
$ws.Range("J16").Value = 0.4476190476190476
$ws.Range("K16").Value = 0.06190476190476191
$ws.Range("M16").Value = 0.01904761904761905
$ws.Range("N16").Value = 0.004761904761904762
$ws.Range("O16").Value = 0.08095238095238096
$ws.Range("S16").Value = 0.1428571428571428
$ws.Range("F17").Value = 0.01162790697674419
$ws.Range("H17").Value = 0.2
$ws.Range("I17").Value = 0.08604651162790698
$ws.Range("J17").Value = 0.4488372093023256
$ws.Range("K17").Value = 0.07209302325581396
$ws.Range("M17").Value = 0.02093023255813953
$ws.Range("N17").Value = 0.002325581395348837
$ws.Range("O17").Value = 0.06744186046511629
$ws.Range("S17").Value = 0.09069767441860466
$ws.Range("F18").Value = 0.02259887005649718
$ws.Range("H18").Value = 0.1581920903954802
$ws.Range("I18").Value = 0.07909604519774012
$ws.Range("J18").Value = 0.4915254237288136
$ws.Range("K18").Value = 0.06779661016949153
$ws.Range("M18").Value = 0.01694915254237288
$ws.Range("O18").Value = 0.05649717514124294
$ws.Range("S18").Value = 0.1073446327683616
$ws.Range("F19").Value = 0.01935483870967742
$ws.Range("H19").Value = 0.2387096774193548
$ws.Range("I19").Value = 0.07004608294930875
$ws.Range("J19").Value = 0.3852534562211982
$ws.Range("K19").Value = 0.09493087557603687
$ws.Range("M19").Value = 0.01658986175115208
$ws.Range("N19").Value = 0.00184331797235023
$ws.Range("O19").Value = 0.06728110599078341
$ws.Range("S19").Value = 0.1059907834101382
